$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 204-205, pushing the existing data (previously
# rows 204-292) down to rows 206-294.
$ws.Rows("204:205").Insert()

# New row 204: Betarraga, Primera quality, date 2021-11-10 (serial 44510)
$ws.Range("A204").Value = 9
$ws.Range("B204").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C204").Value = "Metropolitana"
$ws.Range("D204").Value = 44510
$ws.Range("E204").Value = 13
$ws.Range("F204").Value = 100114014
$ws.Range("G204").Value = "Betarraga"
$ws.Range("H204").Value = "Sin especificar"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 4300
$ws.Range("K204").Value = 90
$ws.Range("L204").Value = 100
$ws.Range("M204").Value = 95
$ws.Range("N204").Value = "`$/unidad"
$ws.Range("O204").Value = "Región Metropolitana"
$ws.Range("P204").Value = 95
$ws.Range("Q204").Value = 1
$ws.Range("R204").Value = "Hortaliza"

# New row 205: Betarraga, Segunda quality, same date
$ws.Range("A205").Value = 9
$ws.Range("B205").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C205").Value = "Metropolitana"
$ws.Range("D205").Value = 44510
$ws.Range("E205").Value = 13
$ws.Range("F205").Value = 100114014
$ws.Range("G205").Value = "Betarraga"
$ws.Range("H205").Value = "Sin especificar"
$ws.Range("I205").Value = "Segunda"
$ws.Range("J205").Value = 2500
$ws.Range("K205").Value = 60
$ws.Range("L205").Value = 70
$ws.Range("M205").Value = 65
$ws.Range("N205").Value = "`$/unidad"
$ws.Range("O205").Value = "Región Metropolitana"
$ws.Range("P205").Value = 65
$ws.Range("Q205").Value = 1
$ws.Range("R205").Value = "Hortaliza"
